$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1056.25
$ws.Range("N40").Value = -1351.5
$ws.Range("L40").Value = 1001.5
$ws.Range("J40").Value = 1001.5
$ws.Range("H43").Value = 1480.2
$ws.Range("I43").Value = 999.75
$ws.Range("M43").Value = -930.75
$ws.Range("K43").Value = 999.75
$ws.Range("L51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("H51").Value = 12400
$ws.Range("N51").ClearContents()
$ws.Range("I55").Value = 172
$ws.Range("J55").Value = 399.33334
$ws.Range("K55").Value = 172
$ws.Range("H55").Value = 269.42856
$ws.Range("N55").Value = -827.33334
$ws.Range("L55").Value = 399.33334
$ws.Range("M55").Value = 42
$ws.Range("I58").Value = 75
$ws.Range("J58").Value = 2429.1428
$ws.Range("K58").Value = 225
$ws.Range("H58").Value = 1906
$ws.Range("N58").Value = -7587.428400000001
$ws.Range("L58").Value = 7287.428400000001
$ws.Range("M58").Value = -75
$ws.Range("H100").Value = 1613.6666
$ws.Range("M100").Value = -629.5
$ws.Range("I100").Value = 1170.5
$ws.Range("K100").Value = 1170.5
$ws.Range("J112").Value = 2000
$ws.Range("H112").Value = 2000
$ws.Range("N112").Value = -8216
$ws.Range("L112").Value = 6000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -8703.444
$ws.Range("K32").Value = 8990.444
$ws.Range("I32").Value = 8990.444
$ws.Range("H32").Value = 12576.167
$ws.Range("M45").Value = -873
$ws.Range("K45").Value = 1250
$ws.Range("I45").Value = 1250
$ws.Range("H45").Value = 1833.3334
$ws.Range("M132").Value = -5914.000100000001
$ws.Range("K132").Value = 8444.000100000001
$ws.Range("I132").Value = 2814.6667
$ws.Range("H132").Value = 4861
$ws.Range("L141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 290
$ws.Range("N10").Value = -478
$ws.Range("K10").Value = 308
$ws.Range("M10").Value = -169
$ws.Range("L10").Value = 200
$ws.Range("I10").Value = 308
$ws.Range("J10").Value = 200
$ws.Range("H33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("I33").Value = 0
$ws.Range("L94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("H94").Value = 3808.3333
$ws.Range("N94").Value = -1902
$ws.Range("J108").Value = 40000
$ws.Range("H108").Value = 40000
$ws.Range("N108").Value = -47680
$ws.Range("L108").Value = 40000
$ws.Range("L121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L113").Value = 5988
$ws.Range("M113").Value = 390.0001
$ws.Range("K113").Value = 1779.9999
$ws.Range("I113").Value = 593.3333
$ws.Range("J113").Value = 1996
$ws.Range("H113").Value = 1294.6666
$ws.Range("N113").Value = -10328
$ws.Range("N132").ClearContents()
$ws.Range("L132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("H132").Value = 900
$ws.Range("I133").Value = 3512
$ws.Range("K133").Value = 10536
$ws.Range("H133").Value = 3512
$ws.Range("M133").Value = -5476

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J80").Value = 27850
$ws.Range("H80").Value = 18794.5
$ws.Range("N80").Value = -29846
$ws.Range("L80").Value = 27850
$ws.Range("L83").Value = 139250
$ws.Range("J83").Value = 27850
$ws.Range("H83").Value = 18794.5
$ws.Range("N83").Value = -149234
$ws.Range("J92").Value = 9985.5
$ws.Range("H92").Value = 9985.5
$ws.Range("N92").Value = -13729.5
$ws.Range("L92").Value = 9985.5
$ws.Range("I122").Value = 0
$ws.Range("H122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("I126").Value = 3481
$ws.Range("M126").Value = -7973
$ws.Range("K126").Value = 10443
$ws.Range("H126").Value = 3481
$ws.Range("N136").Value = -110100
$ws.Range("L136").Value = 105000
$ws.Range("J136").Value = 35000
$ws.Range("H136").Value = 42574

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L16").Value = 1489.8
$ws.Range("K16").Value = 1839.5454
$ws.Range("M16").Value = -1669.5454
$ws.Range("I16").Value = 1839.5454
$ws.Range("J16").Value = 1489.8
$ws.Range("H16").Value = 1730.25
$ws.Range("N16").Value = -1829.8
$ws.Range("N22").Value = -1156.3333
$ws.Range("L22").Value = 566.3333
$ws.Range("M22").Value = -243.25
$ws.Range("I22").Value = 538.25
$ws.Range("J22").Value = 566.3333
$ws.Range("K22").Value = 538.25
$ws.Range("H22").Value = 550.2857
$ws.Range("I27").Value = 538.25
$ws.Range("J27").Value = 566.3333
$ws.Range("K27").Value = 538.25
$ws.Range("H27").Value = 550.2857
$ws.Range("N27").Value = -780.3333
$ws.Range("M27").Value = -431.25
$ws.Range("L27").Value = 566.3333
$ws.Range("I55").Value = 817.5454999999999
$ws.Range("K55").Value = 817.5454999999999
$ws.Range("H55").Value = 841.6429000000001
$ws.Range("M55").Value = -644.5454999999999
$ws.Range("J92").Value = 29000
$ws.Range("H92").Value = 29000
$ws.Range("N92").Value = -33992
$ws.Range("L92").Value = 29000
$ws.Range("N132").ClearContents()
$ws.Range("M132").Value = -15432.5
$ws.Range("L132").Value = 0
$ws.Range("K132").Value = 17962.5
$ws.Range("J132").Value = 0
$ws.Range("I132").Value = 5987.5
$ws.Range("H132").Value = 5987.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M62").Value = -3876
$ws.Range("I62").Value = 4500
$ws.Range("K62").Value = 4500
$ws.Range("H62").Value = 5203.3335
$ws.Range("I65").Value = 4500
$ws.Range("H65").Value = 5203.3335
$ws.Range("M65").Value = -19380
$ws.Range("K65").Value = 22500
$ws.Range("L98").Value = 19295
$ws.Range("J98").Value = 19295
$ws.Range("H98").Value = 19295
$ws.Range("N98").Value = -25285
$ws.Range("H101").Value = 266666
$ws.Range("N101").Value = -273156
$ws.Range("L101").Value = 266666
$ws.Range("J101").Value = 266666
$ws.Range("H109").Value = 187000
$ws.Range("N109").Value = -189774
$ws.Range("L109").Value = 187000
$ws.Range("J109").Value = 187000
$ws.Range("L113").Value = 2250
$ws.Range("J113").Value = 750
$ws.Range("H113").Value = 576
$ws.Range("N113").Value = -6590
$ws.Range("N132").Value = -11357
$ws.Range("M132").Value = -2885.6
$ws.Range("L132").Value = 6297
$ws.Range("K132").Value = 5415.6
$ws.Range("J132").Value = 2099
$ws.Range("I132").Value = 1805.2
$ws.Range("H132").Value = 1854.1666
$ws.Range("N135").Value = -50140
$ws.Range("L135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("H135").Value = 37499.5
